$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has an autofilter / frozen header row and the new data is a
# straightforward extension of the existing table (one more day of Strava
# "weekly scoreboard" rows, 252-259). Copy the last existing row down as a
# template first so number formats (esp. the date style in column B) carry
# over automatically, then overwrite every cell with the real values.
$ws.Range("A251:M251").Copy($ws.Range("A252:M259"))

$ws.Range("A252").Value = "Steven"
$ws.Range("B252").Value = 45494
$ws.Range("C252").Value = "Walk"
$ws.Range("D252").Value = 24
$ws.Range("E252").Value = 1.09
$ws.Range("F252").Value = 46
$ws.Range("G252").Value = 24
$ws.Range("H252").Value = 0
$ws.Range("I252").Value = 0
$ws.Range("J252").Value = 0
$ws.Range("K252").Value = 0
$ws.Range("L252").Value = "Brave Leopard"
$ws.Range("M252").Value = 6

$ws.Range("A253").Value = "Steven"
$ws.Range("B253").Value = 45494
$ws.Range("C253").Value = "Workout"
$ws.Range("D253").Value = 57
$ws.Range("E253").Value = 0
$ws.Range("F253").Value = 0
$ws.Range("G253").Value = 47
$ws.Range("H253").Value = 10
$ws.Range("I253").Value = 0
$ws.Range("J253").Value = 0
$ws.Range("K253").Value = 0
$ws.Range("L253").Value = "Brave Leopard"
$ws.Range("M253").Value = 6

$ws.Range("A254").Value = "Matt"
$ws.Range("B254").Value = 45494
$ws.Range("C254").Value = "Ride"
$ws.Range("D254").Value = 45
$ws.Range("E254").Value = 7.91
$ws.Range("F254").Value = 256
$ws.Range("G254").Value = 19
$ws.Range("H254").Value = 22
$ws.Range("I254").Value = 0
$ws.Range("J254").Value = 0
$ws.Range("K254").Value = 0
$ws.Range("L254").Value = "Wily Hyena"
$ws.Range("M254").Value = 6

$ws.Range("A255").Value = "Matt"
$ws.Range("B255").Value = 45494
$ws.Range("C255").Value = "Ride"
$ws.Range("D255").Value = 45
$ws.Range("E255").Value = 7.92
$ws.Range("F255").Value = 338
$ws.Range("G255").Value = 10
$ws.Range("H255").Value = 28
$ws.Range("I255").Value = 3
$ws.Range("J255").Value = 0
$ws.Range("K255").Value = 0
$ws.Range("L255").Value = "Wily Hyena"
$ws.Range("M255").Value = 6

$ws.Range("A256").Value = "Phil"
$ws.Range("B256").Value = 45494
$ws.Range("C256").Value = "Workout"
$ws.Range("D256").Value = 34
$ws.Range("E256").Value = 0
$ws.Range("F256").Value = 0
$ws.Range("G256").Value = 33
$ws.Range("H256").Value = 1
$ws.Range("I256").Value = 0
$ws.Range("J256").Value = 0
$ws.Range("K256").Value = 0
$ws.Range("L256").Value = "Sauntering Hippo"
$ws.Range("M256").Value = 6

$ws.Range("A257").Value = "Eric"
$ws.Range("B257").Value = 45494
$ws.Range("C257").Value = "Workout"
$ws.Range("D257").Value = 82
$ws.Range("E257").Value = 0
$ws.Range("F257").Value = 0
$ws.Range("G257").Value = 52
$ws.Range("H257").Value = 30
$ws.Range("I257").Value = 1
$ws.Range("J257").Value = 0
$ws.Range("K257").Value = 0
$ws.Range("L257").Value = "Wily Hyena"
$ws.Range("M257").Value = 6

$ws.Range("A258").Value = "Phil"
$ws.Range("B258").Value = 45494
$ws.Range("C258").Value = "Run"
$ws.Range("D258").Value = 31
$ws.Range("E258").Value = 3.1
$ws.Range("F258").Value = 297
$ws.Range("G258").Value = 1
$ws.Range("H258").Value = 9
$ws.Range("I258").Value = 11
$ws.Range("J258").Value = 8
$ws.Range("K258").Value = 0
$ws.Range("L258").Value = "Sauntering Hippo"
$ws.Range("M258").Value = 6

$ws.Range("A259").Value = "Steven"
$ws.Range("B259").Value = 45494
$ws.Range("C259").Value = "Walk"
$ws.Range("D259").Value = 26
$ws.Range("E259").Value = 1.19
$ws.Range("F259").Value = 36
$ws.Range("G259").Value = 26
$ws.Range("H259").Value = 0
$ws.Range("I259").Value = 0
$ws.Range("J259").Value = 0
$ws.Range("K259").Value = 0
$ws.Range("L259").Value = "Brave Leopard"
$ws.Range("M259").Value = 6

# Match the author's final view state: header row frozen, scrolled so the
# newly-appended rows are in view, with the next empty row selected.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A260").Select()